$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 30.39114566666667
$ws.Range("H2").Value = 91.17343700000001
$ws.Range("I2").Value = 0.2485034818803364
$ws.Range("J2").Value = 0.2485034818803363
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.239548
$ws.Range("N2").Value = 0.7186440000000001
$ws.Range("O2").Value = 0.2072176292017679
$ws.Range("P2").Value = 0.2072176292017679
$ws.Range("Q2").Value = 7.280138162158668
$ws.Range("R2").Value = 65.52124345942801
$ws.Range("S2").Value = 0.05149430236362779
$ws.Range("T2").Value = 0.05149430236362779
$ws.Range("G3").Value = 30.39114566666667
$ws.Range("H3").Value = 91.17343700000001
$ws.Range("I3").Value = 0.2485034818803364
$ws.Range("J3").Value = 0.2485034818803363
$ws.Range("O3").Value = 0.1375767575223525
$ws.Range("P3").Value = 0.1375767575223525
$ws.Range("Q3").Value = 4.833458458736112
$ws.Range("R3").Value = 43.501126128625
$ws.Range("S3").Value = 0.03418830327011136
$ws.Range("T3").Value = 0.03418830327011136
$ws.Range("G4").Value = 30.39114566666667
$ws.Range("H4").Value = 91.17343700000001
$ws.Range("I4").Value = 0.2485034818803364
$ws.Range("J4").Value = 0.2485034818803363
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5662826666666666
$ws.Range("N4").Value = 1.698848
$ws.Range("O4").Value = 0.4898548585031879
$ws.Range("P4").Value = 0.4898548585031879
$ws.Range("Q4").Value = 17.20997901117511
$ws.Range("R4").Value = 154.889811100576
$ws.Range("S4").Value = 0.1217306379540417
$ws.Range("T4").Value = 0.1217306379540417
$ws.Range("G5").Value = 30.39114566666667
$ws.Range("H5").Value = 91.17343700000001
$ws.Range("I5").Value = 0.2485034818803364
$ws.Range("J5").Value = 0.2485034818803363
$ws.Range("M5").Value = 0.191149
$ws.Range("N5").Value = 0.573447
$ws.Range("O5").Value = 0.1653507547726916
$ws.Range("P5").Value = 0.1653507547726916
$ws.Range("Q5").Value = 5.809237103037668
$ws.Range("R5").Value = 52.28313392733901
$ws.Range("S5").Value = 0.04109023829255552
$ws.Range("T5").Value = 0.04109023829255551
$ws.Range("G6").Value = 45.91529066666667
$ws.Range("I6").Value = 0.3754419042757282
$ws.Range("J6").Value = 0.3754419042757282
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.239548
$ws.Range("N6").Value = 0.7186440000000001
$ws.Range("O6").Value = 0.2072176292017679
$ws.Range("P6").Value = 0.2072176292017679
$ws.Range("Q6").Value = 10.99891604861867
$ws.Range("R6").Value = 98.99024443756802
$ws.Range("S6").Value = 0.07779818130701348
$ws.Range("T6").Value = 0.07779818130701348
$ws.Range("G7").Value = 45.91529066666667
$ws.Range("I7").Value = 0.3754419042757282
$ws.Range("J7").Value = 0.3754419042757282
$ws.Range("O7").Value = 0.1375767575223525
$ws.Range("P7").Value = 0.1375767575223525
$ws.Range("R7").Value = 65.72199917800002
$ws.Range("S7").Value = 0.05165207982827213
$ws.Range("T7").Value = 0.05165207982827215
$ws.Range("G8").Value = 45.91529066666667
$ws.Range("I8").Value = 0.3754419042757282
$ws.Range("J8").Value = 0.3754419042757282
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5662826666666666
$ws.Range("N8").Value = 1.698848
$ws.Range("O8").Value = 0.4898548585031879
$ws.Range("P8").Value = 0.4898548585031879
$ws.Range("Q8").Value = 26.00103323949511
$ws.Range("R8").Value = 234.009299155456
$ws.Range("S8").Value = 0.1839120408951542
$ws.Range("T8").Value = 0.1839120408951543
$ws.Range("G9").Value = 45.91529066666667
$ws.Range("I9").Value = 0.3754419042757282
$ws.Range("J9").Value = 0.3754419042757282
$ws.Range("M9").Value = 0.191149
$ws.Range("N9").Value = 0.573447
$ws.Range("O9").Value = 0.1653507547726916
$ws.Range("P9").Value = 0.1653507547726916
$ws.Range("Q9").Value = 8.776661895642668
$ws.Range("R9").Value = 78.98995706078402
$ws.Range("S9").Value = 0.06207960224528829
$ws.Range("T9").Value = 0.06207960224528829
$ws.Range("G10").Value = 2.332475
$ws.Range("H10").Value = 6.997425
$ws.Range("I10").Value = 0.01907227076123622
$ws.Range("J10").Value = 0.01907227076123622
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.239548
$ws.Range("N10").Value = 0.7186440000000001
$ws.Range("O10").Value = 0.2072176292017679
$ws.Range("P10").Value = 0.2072176292017679
$ws.Range("Q10").Value = 0.5587397213
$ws.Range("R10").Value = 5.028657491700001
$ws.Range("S10").Value = 0.003952110730637567
$ws.Range("T10").Value = 0.003952110730637567
$ws.Range("G11").Value = 2.332475
$ws.Range("H11").Value = 6.997425
$ws.Range("I11").Value = 0.01907227076123622
$ws.Range("J11").Value = 0.01907227076123622
$ws.Range("O11").Value = 0.1375767575223525
$ws.Range("P11").Value = 0.1375767575223525
$ws.Range("Q11").Value = 0.3709607114583333
$ws.Range("R11").Value = 3.338646403125
$ws.Range("S11").Value = 0.002623901169919249
$ws.Range("T11").Value = 0.002623901169919249
$ws.Range("G12").Value = 2.332475
$ws.Range("H12").Value = 6.997425
$ws.Range("I12").Value = 0.01907227076123622
$ws.Range("J12").Value = 0.01907227076123622
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.5662826666666666
$ws.Range("N12").Value = 1.698848
$ws.Range("O12").Value = 0.4898548585031879
$ws.Range("P12").Value = 0.4898548585031879
$ws.Range("Q12").Value = 1.320840162933333
$ws.Range("R12").Value = 11.8875614664
$ws.Range("S12").Value = 0.009342644495079856
$ws.Range("T12").Value = 0.009342644495079856
$ws.Range("G13").Value = 2.332475
$ws.Range("H13").Value = 6.997425
$ws.Range("I13").Value = 0.01907227076123622
$ws.Range("J13").Value = 0.01907227076123622
$ws.Range("M13").Value = 0.191149
$ws.Range("N13").Value = 0.573447
$ws.Range("O13").Value = 0.1653507547726916
$ws.Range("P13").Value = 0.1653507547726916
$ws.Range("Q13").Value = 0.445850263775
$ws.Range("R13").Value = 4.012652373975
$ws.Range("S13").Value = 0.003153614365599547
$ws.Range("T13").Value = 0.003153614365599547
$ws.Range("G14").Value = 43.657748
$ws.Range("H14").Value = 130.973244
$ws.Range("I14").Value = 0.3569823430826993
$ws.Range("J14").Value = 0.3569823430826993
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.239548
$ws.Range("N14").Value = 0.7186440000000001
$ws.Range("O14").Value = 0.2072176292017679
$ws.Range("P14").Value = 0.2072176292017679
$ws.Range("Q14").Value = 10.458126217904
$ws.Range("R14").Value = 94.123135961136
$ws.Range("S14").Value = 0.07397303480048908
$ws.Range("T14").Value = 0.07397303480048908
$ws.Range("G15").Value = 43.657748
$ws.Range("H15").Value = 130.973244
$ws.Range("I15").Value = 0.3569823430826993
$ws.Range("J15").Value = 0.3569823430826993
$ws.Range("O15").Value = 0.1375767575223525
$ws.Range("P15").Value = 0.1375767575223525
$ws.Range("Q15").Value = 6.943401004833333
$ws.Range("R15").Value = 62.4906090435
$ws.Range("S15").Value = 0.04911247325404977
$ws.Range("T15").Value = 0.04911247325404978
$ws.Range("G16").Value = 43.657748
$ws.Range("H16").Value = 130.973244
$ws.Range("I16").Value = 0.3569823430826993
$ws.Range("J16").Value = 0.3569823430826993
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.5662826666666666
$ws.Range("N16").Value = 1.698848
$ws.Range("O16").Value = 0.4898548585031879
$ws.Range("P16").Value = 0.4898548585031879
$ws.Range("Q16").Value = 24.72262595810133
$ws.Range("R16").Value = 222.503633622912
$ws.Range("S16").Value = 0.1748695351589121
$ws.Range("T16").Value = 0.1748695351589122
$ws.Range("G17").Value = 43.657748
$ws.Range("H17").Value = 130.973244
$ws.Range("I17").Value = 0.3569823430826993
$ws.Range("J17").Value = 0.3569823430826993
$ws.Range("M17").Value = 0.191149
$ws.Range("N17").Value = 0.573447
$ws.Range("O17").Value = 0.1653507547726916
$ws.Range("P17").Value = 0.1653507547726916
$ws.Range("Q17").Value = 8.345134872452
$ws.Range("R17").Value = 75.106213852068
$ws.Range("S17").Value = 0.05902729986924828
$ws.Range("T17").Value = 0.05902729986924828
